$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 34: Kenya COVID-19 update for 16th April 2020 (date serial 43937)
$ws.Range("A34").Value = 43937
$ws.Range("A34").NumberFormat = "d-mmm-yy"
$ws.Range("B34").Value = 9
$ws.Range("C34").Value = 704
$ws.Range("F34").Value = 234
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("L34").Value = "2-64."

# Update selection / view to match the saved state after the edit
$ws.Range("E33").Select()
$excel.ActiveWindow.ScrollRow = 18
$excel.ActiveWindow.ScrollColumn = 1
